$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the two price figures in column B that were mis-entered.
$ws.Range("B3").Value = '$91.00 - $252.75'
$ws.Range("B4").Value = ' 73 00-'

# Leave the selection on the first cell instead of the stale B1:F6 block.
[void]$ws.Range("A1").Select()
